$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.159.17"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").Value = "3.804.52"
$ws.Range("E3").Value = "  -1.02%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "701.10"
$ws.Range("E5").Value = "  +0.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.14"
$ws.Range("E6").Value = "  -1.44%  "
$ws.Range("D7").Value = "3.806.21"
$ws.Range("E7").Value = "  -0.79%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.524"
$ws.Range("E9").Value = "  -0.46%  "
$ws.Range("E10").Value = "  -1.39%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.54"
$ws.Range("E11").Value = "  +2.95%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.478"
$ws.Range("E12").Value = "  +4.06%  "
$ws.Range("E13").Value = "  -1.93%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.30"
$ws.Range("E14").Value = "  -1.43%  "
$ws.Range("D15").Value = "4.445.45"
$ws.Range("E15").Value = "  -1.16%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "71.413.56"
$ws.Range("E16").Value = "  +0.37%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.688.02"
$ws.Range("E17").Value = "  -3.86%  "
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.60"
$ws.Range("E18").Value = "  +0.92%  "
$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.20"
$ws.Range("E19").Value = "  -0.53%  "
$ws.Range("E20").Value = "  +0.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "514.64"
$ws.Range("E21").Value = "  +3.34%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.44"
$ws.Range("E22").Value = "  -2.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.715"
$ws.Range("E23").Value = "  -2.89%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.71"
$ws.Range("E24").Value = "  -1.89%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000141"
$ws.Range("E25").Value = "  -2.87%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.59"
$ws.Range("E26").Value = "  +3.40%  "
$ws.Range("D27").Value = "3.949.00"
$ws.Range("E27").Value = "  -1.34%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.25"
$ws.Range("E28").Value = "  -3.61%  "
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.99"
$ws.Range("E30").Value = "  -5.11%  "
$ws.Range("E31").Value = "  -4.00%  "
$ws.Range("E32").Value = "  +1.18%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.29"
$ws.Range("E33").Value = "  -2.34%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.22"
$ws.Range("E34").Value = "  -0.61%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.171"
$ws.Range("E35").Value = "  -3.39%  "
$ws.Range("E36").Value = "  +1.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.21%  "
$ws.Range("D38").Value = "3.766.52"
$ws.Range("E38").Value = "  -1.19%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.69"
$ws.Range("E39").Value = "  +11.68%  "
$ws.Range("E40").Value = "  -2.20%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.36"
$ws.Range("E41").Value = "  +1.16%  "
$ws.Range("E42").Value = "  -2.43%  "
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.18"
$ws.Range("E44").Value = "  -4.87%  "
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "163.75"
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "49.35"
$ws.Range("E47").Value = "  +0.72%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.000303"
$ws.Range("E48").Value = "  -4.19%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "419.72"
$ws.Range("E49").Value = "  -3.13%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.39"
$ws.Range("E50").Value = "  +0.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.65"
$ws.Range("E51").Value = "  -1.06%  "
